# Remove the extra "Valor do ganho do sensor variando em ± 20%" banner row
# that sits above the "Obs.: ..." note on the "Melhor rede" sheet. Deleting
# the entire row shifts every row below it up by one, matching the
# target layout (Obs. row becomes row 2, "Melhores Redes Treinadas" becomes
# row 4, the data table becomes rows 5-8, and the summary block becomes
# rows 10-12).
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Melhor rede")
if ($ws -eq $null) {
    $ws = $wb.ActiveSheet
}

$ws.Rows.Item(2).Delete()

# The color-scale conditional format on the results table tracked the old
# H7:H9 range; re-anchor it to the new H6:H8 location now that the rows
# have shifted up.
$newRange = $ws.Range("H6:H8")
$fcs = $newRange.FormatConditions
if ($fcs.Count -ge 1) {
    $fc = $fcs.Item(1)
    $fc.ModifyAppliesToRange($newRange)
}
